$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final touches before hand-in:
# - Row 10 (test case "Insert buy" placeholder) now documents the
#   "Retrieve price from control layer" test case instead.
$ws.Range("B10").Value = "Retrieve price from control layer"

# The row grew a touch taller to fit the new wrapped text.
$ws.Rows.Item(10).RowHeight = 48

# Leave the cursor parked below the table, as it was when the author saved.
$ws.Range("B21").Select() | Out-Null
